$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "300.90", "42.620.83") that must
# stay as literal text, matching the source inlineStr cells. A leading
# apostrophe is Excel standard "force text" entry convention: it keeps the
# cell text-typed with a General number format (just flags quotePrefix),
# instead of letting Excel silently parse the string into a numeric value.

$ws.Range('D2').Value = '''42.620.83'
$ws.Range('E2').Value = '  -1.57%  '

$ws.Range('D3').Value = '''2.284.72'
$ws.Range('E3').Value = '  -3.44%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '''300.90'
$ws.Range('E5').Value = '  -2.85%  '

$ws.Range('D6').Value = '''97.12'
$ws.Range('E6').Value = '  -6.33%  '

$ws.Range('E7').Value = '  -1.81%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '''0.501'
$ws.Range('E9').Value = '  -3.65%  '

$ws.Range('D10').Value = '''33.50'
$ws.Range('E10').Value = '  -5.95%  '

$ws.Range('D11').Value = '''0.0787'
$ws.Range('E11').Value = '  -2.11%  '

$ws.Range('D12').Value = '''50.57'
$ws.Range('E12').Value = '  -4.85%  '

$ws.Range('E13').Value = '  -0.12%  '

$ws.Range('D14').Value = '''6.64'
$ws.Range('E14').Value = '  -3.99%  '

$ws.Range('D15').Value = '''2.641.43'
$ws.Range('E15').Value = '  -3.52%  '

$ws.Range('D16').Value = '''15.20'
$ws.Range('E16').Value = '  -1.87%  '

$ws.Range('D17').Value = '''2.298.31'
$ws.Range('E17').Value = '  -2.88%  '

$ws.Range('D18').Value = '''0.788'
$ws.Range('E18').Value = '  -2.50%  '

$ws.Range('D19').Value = '''42.554.77'
$ws.Range('E19').Value = '  -1.70%  '

$ws.Range('D20').Value = '''0.0₃0895'

$ws.Range('D21').Value = '''11.49'
$ws.Range('E21').Value = '  -3.20%  '

$ws.Range('E22').Value = '  -5.11%  '

$ws.Range('E23').Value = '  -1.98%  '

$ws.Range('D24').Value = '''235.11'
$ws.Range('E24').Value = '  -1.92%  '

$ws.Range('E25').Value = '  -4.67%  '

$ws.Range('E26').Value = '  -4.15%  '

$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('D28').Value = '''24.48'
$ws.Range('E28').Value = '  -4.97%  '

$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = '''165.84'
$ws.Range('E29').Value = '  +2.69%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.07'
$ws.Range('E30').Value = '  -11.11%  '

$ws.Range('D31').Value = '''33.66'
$ws.Range('E31').Value = '  -7.93%  '

$ws.Range('D32').Value = '''9.10'
$ws.Range('E32').Value = '  -3.63%  '

$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('D34').Value = '''4.97'
$ws.Range('E34').Value = '  -4.52%  '

$ws.Range('E35').Value = '  -3.87%  '

$ws.Range('D36').Value = '''0.0695'
$ws.Range('E36').Value = '  -5.32%  '

$ws.Range('D37').Value = '''4.36'
$ws.Range('E37').Value = '  -6.42%  '

$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''2.83'
$ws.Range('E38').Value = '  -8.00%  '

$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '''16.24'
$ws.Range('E39').Value = '  -10.36%  '

$ws.Range('E40').Value = '  -7.40%  '

$ws.Range('D41').Value = '''0.100'
$ws.Range('E41').Value = '  -4.68%  '

$ws.Range('E42').Value = '  -2.90%  '

$ws.Range('D43').Value = '''2.42'
$ws.Range('E43').Value = '  -6.55%  '

$ws.Range('D44').Value = '''1.962.91'
$ws.Range('E44').Value = '  -3.75%  '

$ws.Range('E45').Value = '  -2.21%  '

$ws.Range('D46').Value = '''17.90'
$ws.Range('E46').Value = '  -8.03%  '

$ws.Range('D47').Value = '''9.67'
$ws.Range('E47').Value = '  -8.56%  '

$ws.Range('D48').Value = '''2.83'
$ws.Range('E48').Value = '  -8.07%  '

$ws.Range('D49').Value = '''53.27'
$ws.Range('E49').Value = '  -7.65%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '''2.507.98'
$ws.Range('E51').Value = '  -3.46%  '

Write-Output "Applied cryptos.xlsx update"
